$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Locate the paragraphs we need to touch by their text content so the
# script does not depend on a brittle, hard-coded paragraph index.
# -----------------------------------------------------------------
$titleIndex = 0
$logoIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd()
    if ($t -eq "QuIz de REACT JS") {
        $titleIndex = $i
    }
    if ($t -eq "Agregar logo al encabezado") {
        $logoIndex = $i
    }
}

# -----------------------------------------------------------------
# 1) Title paragraph: "QuIz de REACT JS" -> "Quiz de REACT JS", split
#    across three runs ("Qu" / "i" / "z de REACT JS") and drop the
#    spell-check proofErr markers that used to wrap "QuIz".
# -----------------------------------------------------------------
$p1 = $d.Paragraphs($titleIndex)
$fullTitle = $d.Range($p1.Range.Start, $p1.Range.End)
$titleStart = $p1.Range.Start
$fullTitle.Delete()

$d.Paragraphs($titleIndex).Range.InsertParagraphBefore()
$newTitle = $d.Paragraphs($titleIndex)
$newTitle.Range.InsertBefore("Quiz de REACT JS")

# Force a run boundary around the lowercase "i" so it ends up as its
# own <w:r>, matching "Qu" | "i" | "z de REACT JS".
$iRange = $d.Range($titleStart + 2, $titleStart + 3)
$iRange.Bold = $true
$iRange.Bold = $false

# -----------------------------------------------------------------
# 2) Split "Agregar logo al encabezado" list item: the trailing
#    _GoBack bookmark moves into a brand-new list paragraph that
#    reads "Agregar ".
# -----------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$agregarLogo = $d.Paragraphs($logoIndex)
$agregarLogo.Range.InsertParagraphAfter()
$newItem = $d.Paragraphs($logoIndex + 1)

# Insert a placeholder char after "Agregar " so the bookmark can be
# planted at a safe mid-text offset, then trim the placeholder back
# off so the bookmark ends up collapsed right before the paragraph
# mark (mirrors the original "Agregar logo al encabezado" + _GoBack
# layout).
$newItem.Range.InsertBefore("Agregar X")
$newItem2 = $d.Paragraphs($logoIndex + 1)
$bmPos = $newItem2.Range.Start + 8
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$placeholder = $d.Range($bmPos, $bmPos + 1)
$placeholder.Delete()
